# Add data for 2024-12-13 (new column CR) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column CR (column index 96) ---
# Match the width of the preceding data columns (e.g. CQ / column 95).
$ws.Cells.Item(1, 96).EntireColumn.ColumnWidth = $ws.Cells.Item(1, 95).EntireColumn.ColumnWidth

# --- Header cell CR1: "2024/12/13" as text (not an auto-converted date) ---
$headerDst = $ws.Cells.Item(1, 96)
$headerSrc = $ws.Cells.Item(1, 95)   # CQ1 already has the desired style (s=1)
$headerDst.NumberFormat = "@"
$headerDst.Value = "2024/12/13"
$headerSrc.Copy()
$headerDst.PasteSpecial(-4122)      # xlPasteFormats: copy CQ1's style onto CR1
$excel.CutCopyMode = 0

# --- Style templates taken from row 2, one cell per fill-threshold style ---
# s=1 : value >= 140   (no special fill)
# s=2 : value < 125     (yellow fill)
# s=3 : 125 <= value < 140 (light-blue fill)
$tmplStyle1 = $ws.Cells.Item(2, 1)    # A2  (478, s=1)
$tmplStyle2 = $ws.Cells.Item(2, 4)    # D2  (121.8, s=2)
$tmplStyle3 = $ws.Cells.Item(2, 14)   # N2  (129.8, s=3)

# --- New data for 2024/12/13, rows 2-53 ---
$data = @(
    @{ Row=2; Style=3; Val=134.9 },
    @{ Row=3; Style=3; Val=131.6 },
    @{ Row=4; Style=1; Val=140.6 },
    @{ Row=5; Style=3; Val=131.5 },
    @{ Row=6; Style=1; Val=164.1 },
    @{ Row=7; Style=1; Val=255.3 },
    @{ Row=8; Style=1; Val=156.4 },
    @{ Row=9; Style=3; Val=137.2 },
    @{ Row=10; Style=1; Val=177.5 },
    @{ Row=11; Style=1; Val=173.8 },
    @{ Row=12; Style=1; Val=151.1 },
    @{ Row=13; Style=1; Val=200 },
    @{ Row=14; Style=1; Val=165.1 },
    @{ Row=15; Style=1; Val=146.3 },
    @{ Row=16; Style=3; Val=139.9 },
    @{ Row=17; Style=1; Val=152.5 },
    @{ Row=18; Style=1; Val=167.3 },
    @{ Row=19; Style=1; Val=162.5 },
    @{ Row=20; Style=1; Val=140 },
    @{ Row=21; Style=1; Val=259.2 },
    @{ Row=22; Style=1; Val=171.5 },
    @{ Row=23; Style=1; Val=166.5 },
    @{ Row=24; Style=1; Val=142.3 },
    @{ Row=25; Style=1; Val=214.4 },
    @{ Row=26; Style=1; Val=152.2 },
    @{ Row=27; Style=1; Val=140.9 },
    @{ Row=28; Style=1; Val=464.3 },
    @{ Row=29; Style=1; Val=162 },
    @{ Row=30; Style=3; Val=132 },
    @{ Row=31; Style=1; Val=195.3 },
    @{ Row=32; Style=1; Val=140.3 },
    @{ Row=33; Style=1; Val=201 },
    @{ Row=34; Style=1; Val=152.5 },
    @{ Row=35; Style=1; Val=164.1 },
    @{ Row=36; Style=1; Val=141.9 },
    @{ Row=37; Style=1; Val=194.1 },
    @{ Row=38; Style=3; Val=131.4 },
    @{ Row=39; Style=1; Val=254.6 },
    @{ Row=40; Style=1; Val=159.4 },
    @{ Row=41; Style=1; Val=244.9 },
    @{ Row=42; Style=1; Val=151.9 },
    @{ Row=43; Style=1; Val=175.4 },
    @{ Row=44; Style=1; Val=142.4 },
    @{ Row=45; Style=1; Val=187.5 },
    @{ Row=46; Style=1; Val=261.6 },
    @{ Row=47; Style=2; Val=123.4 },
    @{ Row=48; Style=1; Val=236.4 },
    @{ Row=49; Style=1; Val=202.2 },
    @{ Row=50; Style=1; Val=142 },
    @{ Row=51; Style=2; Val=123.5 },
    @{ Row=52; Style=1; Val=179.5 },
    @{ Row=53; Style=1; Val=140.4 }
)

foreach ($item in $data) {
    $dst = $ws.Cells.Item($item.Row, 96)
    $dst.Value = $item.Val
    if ($item.Style -eq 1) {
        $tmplStyle1.Copy()
    } elseif ($item.Style -eq 2) {
        $tmplStyle2.Copy()
    } else {
        $tmplStyle3.Copy()
    }
    $dst.PasteSpecial(-4122)
}
$excel.CutCopyMode = 0
